# "param study general work"
# Adds 8 new parameter-study rows (67-74, copies of row 2 with the
# sequence index and one varied input updated), turns the new extent
# into an AutoFilter range (with the corresponding hidden
# _FilterDatabase defined name), and updates the current selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues  = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# New rows are all copies of the first data row (row 2), only the CASE
# index (col A) and a single varied parameter change per row.
$newRows = @(
    @{ Row = 67; Case = 66; Col = "B"; Value = 0.005 },
    @{ Row = 68; Case = 67; Col = "B"; Value = 0.015 },
    @{ Row = 69; Case = 68; Col = "B"; Value = 0.158 },
    @{ Row = 70; Case = 69; Col = "B"; Value = 0.5 },
    @{ Row = 71; Case = 70; Col = "D"; Value = 4.2 },
    @{ Row = 72; Case = 71; Col = "D"; Value = 35.5 },
    @{ Row = 73; Case = 72; Col = "D"; Value = 50 },
    @{ Row = 74; Case = 73; Col = "D"; Value = 76.5 }
)

foreach ($item in $newRows) {
    $destRow = $item.Row

    # Clone the whole baseline row (values first, then re-stamp the
    # N:O number-format styling that Copy/PasteAll doesn't carry over).
    $ws.Range("A2:O2").Copy()
    $ws.Range("A$($destRow):O$($destRow)").PasteSpecial($xlPasteValues)
    $ws.Range("N2:O2").Copy()
    $ws.Range("N$($destRow):O$($destRow)").PasteSpecial($xlPasteFormats)

    $ws.Cells.Item($destRow, 1).Value = $item.Case
    $ws.Cells.Item($destRow, [int][char]$item.Col - [int][char]"A" + 1).Value = $item.Value
}

$excel.CutCopyMode = $false

# Extend the AutoFilter over the whole (now larger) table and register
# the resulting hidden _FilterDatabase defined name on the sheet.
$ws.Range("A1:O74").AutoFilter()
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=paramStudy!`$A`$1:`$O`$74")
$filterName.Visible = $false

# Match the author's final selection.
$ws.Range("E71").Select()
